$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = -21.183
$ws.Range("E5").Value = 12.94
$ws.Range("E9").Value = 12.82
$ws.Range("E11").Value = 13.246
$ws.Range("A21").Value = -20.88
$ws.Range("E21").Value = 13.535
$ws.Range("A23").Value = -21.376
$ws.Range("A25").Value = -22.27
